# Apply the "mise à jour avant publication" changes to dist_aide.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the displayed hyperlink text in column C -------------------
# (the underlying hyperlink targets in the worksheet relationships are left
#  untouched - only the visible text of the cells changes)

# Row 8: replace the displayed "menu.pdf" link text with the new help file
$ws.Range("C8").Value2 = "https://ductair.github.io/ductaironline/Support/Aide_transfo_trémie.pdf"

# Row 5: fix displayed "dournisseurs.pdf" -> "fournisseurs.pdf"
$ws.Range("C5").Value2 = "https://ductair.github.io/ductaironline/Support/fournisseurs.pdf"

# --- Refresh the "Date_lien" column (D2:D9) -----------------------------
$ws.Range("D2:D9").Value2 = 45748

# --- Update the active selection left over in the sheet view -----------
$ws.Activate()
$ws.Range("B10").Select()
